$wb = $excel.ActiveWorkbook

# The three "zones" sheets (zones priority / zones economy / zones cp) each have a
# Country / iso / Zone table in columns A:C. The edit removes the "iso" country-code
# column (column B), shifting the "Zone" column from C to B.
$zoneSheetNames = @("zones priority ", "zones economy", "zones cp")
foreach ($name in $zoneSheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Columns.Item(2).EntireColumn.Delete()
}

# Update the lingering selections on each sheet to match the saved view state.
$wsPriorityZones = $wb.Worksheets.Item("zones priority ")
$wsPriorityZones.Activate()
$wsPriorityZones.Range("C24").Select()

$wsEconomyZones = $wb.Worksheets.Item("zones economy")
$wsEconomyZones.Activate()
$wsEconomyZones.Columns.Item(2).Select()

# "zones cp" ends up being the active/selected sheet when the file was saved.
$wsCpZones = $wb.Worksheets.Item("zones cp")
$wsCpZones.Activate()
$wsCpZones.Range("C9").Select()
